$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.447.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").Value = "'1.632.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").Value = "'305.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("D7").Value = "'0.3753"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("D8").Value = "'0.3673"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").Value = "'51.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("D10").Value = "'0.08198"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").Value = "'1.230"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.80%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").Value = "'22.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").Value = "'6.578"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("D15").Value = "'0.00001252"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.64%  "

$ws.Range("D16").Value = "'7.286"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.99%  "

$ws.Range("D17").Value = "'1.635.43"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'94.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("E19").Value = "  +0.80%  "

$ws.Range("D20").Value = "'17.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("D21").Value = "'6.470"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.65%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "'12.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "

$ws.Range("D24").Value = "'23.444.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("D25").Value = "'3.179"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.20%  "

$ws.Range("D26").Value = "'2.463"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("D27").Value = "'21.44"
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = "  -0.70%  "

$ws.Range("D29").Value = "'5.327"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "

$ws.Range("D30").Value = "'134.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").Value = "'1.817.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "

$ws.Range("E32").Value = "  -4.87%  "

$ws.Range("D33").Value = "'6.835"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").Value = "'1.037"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.98%  "

$ws.Range("E35").Value = "  +5.32%  "

$ws.Range("D36").Value = "'0.02794"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.69%  "

$ws.Range("D37").Value = "'0.2538"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").Value = "'0.08791"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "

$ws.Range("D39").Value = "'6.091"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.99%  "

$ws.Range("D40").Value = "'0.07147"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.22%  "

$ws.Range("D41").Value = "'0.7084"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").Value = "'1.351"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.63%  "

$ws.Range("D43").Value = "'16.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("D44").Value = "'12.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.91%  "

$ws.Range("D45").Value = "'0.6570"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").Value = "'2.337"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.58%  "

$ws.Range("D47").Value = "'0.9997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "'3.996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.23%  "

$ws.Range("D49").Value = "'0.08040"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.57%  "

$ws.Range("D50").Value = "'1.211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("D51").Value = "'125.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.43%  "
